$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format before writing so that
# numeric-looking strings (e.g. "1.001") are stored as literal text
# rather than being coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.523.43"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.744.24"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "323.08"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4449"
$ws.Range("E7").Value = "  +4.68%  "
$ws.Range("D8").Value = "0.3521"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "0.07420"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value = "41.57"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").Value = "1.078"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "20.47"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "5.906"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "7.079"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "1.741.74"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "91.50"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "0.00001053"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "0.06384"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "16.84"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "5.717"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").Value = "27.549.56"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").Value = "2.095"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "160.59"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "20.04"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "1.942.25"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "125.15"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "2.030"
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("D31").Value = "1.047"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("D32").Value = "0.09080"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "3.652"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "5.366"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").Value = "0.02273"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("D37").Value = "0.06032"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "0.2063"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "4.893"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "0.6234"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "1.184"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "7.723"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "3.701"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "0.5792"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "121.99"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "1.924"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "1.114"
$ws.Range("E50").Value = "  -4.29%  "
$ws.Range("D51").Value = "71.36"
$ws.Range("E51").Value = "  -2.28%  "

# Restore the default style on the Price column so formatting matches
# the original workbook (only the values themselves changed).
$ws.Range("D2:D51").Style = "Normal"
